$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert A1 header back to "student_name" (was changed to "student_number")
$ws.Range("A1").Value = "student_name"

# Restore the active selection to D4 as it was at last save
$ws.Range("D4").Select()
